# Reorders the "Recorded By" (column G) list of names/emails in the
# Session Analysis Results sheet:
#   - if "dnasr281@gmail.com" is present in the comma-separated list,
#     move it to the front;
#   - if "backup@backdoor.com" is present in the (resulting) list,
#     move it to the end.
# Entries that are a single value, or that do not contain either of
# those two tokens, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $value = $cell.Value2

    if ($null -eq $value -or $value -eq "") {
        continue
    }

    $parts = @()
    foreach ($p in $value -split ",") {
        $parts += $p.Trim()
    }

    if ($parts.Count -le 1) {
        continue
    }

    if ($parts -contains "dnasr281@gmail.com") {
        $rest = @($parts | Where-Object { $_ -ne "dnasr281@gmail.com" })
        $parts = @("dnasr281@gmail.com") + $rest
    }

    if ($parts -contains "backup@backdoor.com") {
        $rest = @($parts | Where-Object { $_ -ne "backup@backdoor.com" })
        $parts = $rest + @("backup@backdoor.com")
    }

    $newValue = $parts -join ", "

    if ($newValue -ne $value) {
        $cell.Value = $newValue
    }
}
